# AFDP-892 - Store container folder IDs in a single shared table, not in
# columns in the container's table - convert CaseFile module to use
# AcmContainerFolder.
#
# Update the two rule-condition cells on the "Save Case File Rules" rule
# table (Sheet1) that referenced the old per-table ECM folder / due date
# columns so they reference the new shared AcmContainerFolder / dueDate
# fields instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Assign Alfresco Folder" rule condition (row 18, column C)
$ws.Range("C18").Value = "containerFolder.cmisFolderId == null"

# "Set Due Date" rule condition (row 21, column C)
$ws.Range("C21").Value = "dueDate == null"

# Row heights drifted slightly as a side effect of editing the sheet in
# Excel; match the committed values.
$ws.Rows.Item(18).RowHeight = 13.8
$ws.Rows.Item(21).RowHeight = 13.8

# Selection / scroll position left by the editor.
[void]$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$win.TabRatio = 0.185
[void]$ws.Range("D20").Select()
